# Sprint 2 backlog update: mark tasks 4 ("User can edit/delete comment"),
# 5 ("User can create a room"), 8 ("User can vote for a post") and
# 9 ("User can vote for a comment") with their new status, carrying over
# the matching conditional-style formatting (fill/border/font) used
# elsewhere in the sheet for each status color.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 8: restyle to the "green" status block, set status text to InProgress ---
$ws.Range("B6:F6").Copy() | Out-Null
$ws.Range("B8:F8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F8").Value = "InProgress"

# --- Row 9: restyle to the "orange" status block, set status text to InProgress ---
$ws.Range("B17:F17").Copy() | Out-Null
$ws.Range("B9:F9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F9").Value = "InProgress"

# --- Row 12: restyle to the "green" status block, set status text to Complete ---
$ws.Range("B6:F6").Copy() | Out-Null
$ws.Range("B12:F12").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("F12").Value = "Complete"

# --- Row 13: restyle to the "orange" status block, set status text to InProgress ---
$ws.Range("B17:F17").Copy() | Out-Null
$ws.Range("B13:F13").PasteSpecial(-4122) | Out-Null # xlPasteFormats
$ws.Range("F13").Value = "InProgress"

$excel.CutCopyMode = 0

# --- Update the active view: scroll near the top and select F10 ---
$ws.Range("F10").Select() | Out-Null

Write-Host "Backlog rows 8, 9, 12 and 13 updated for sprint 2"
